$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.449729800224304
$ws.Range("B1").Value = 4.660916805267334
$ws.Range("C1").Value = 2.269258499145508
$ws.Range("D1").Value = 1.628582954406738
$ws.Range("E1").Value = 1.405364513397217
